$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '99.268.03'
$ws.Range("E2").Value = '  +1.05%  '
$ws.Range("D3").Value = '3.277.94'
$ws.Range("E3").Value = '  -2.39%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '253.32'
$ws.Range("E5").Value = '  -0.67%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '621.59'
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("E7").Value = '  +21.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.399'
$ws.Range("E8").Value = '  +3.36%  '
$ws.Range("E9").Value = '  -0.03%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.966'
$ws.Range("E10").Value = '  +20.60%  '
$ws.Range("D11").Value = '3.274.55'
$ws.Range("E11").Value = '  -2.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.200'
$ws.Range("E12").Value = '  +0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.23'
$ws.Range("E13").Value = '  +9.38%  '
$ws.Range("D14").Value = '98.889.75'
$ws.Range("E14").Value = '  +0.92%  '
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").Value = '3.868.48'
$ws.Range("E16").Value = '  -2.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.46'
$ws.Range("E17").Value = '  -0.92%  '
$ws.Range("D18").Value = '3.252.70'
$ws.Range("E18").Value = '  -3.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.41'
$ws.Range("E19").Value = '  -5.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.25'
$ws.Range("E20").Value = '  +1.83%  '
$ws.Range("E21").Value = '  +7.75%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '486.84'
$ws.Range("E22").Value = '  +0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.27'
$ws.Range("E23").Value = '  +1.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000200'
$ws.Range("E24").Value = '  -3.21%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.62'
$ws.Range("E25").Value = '  -0.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '88.92'
$ws.Range("E26").Value = '  +0.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.320'
$ws.Range("E27").Value = '  +29.97%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.94'
$ws.Range("E28").Value = '  -0.72%  '
$ws.Range("D29").Value = '3.424.62'
$ws.Range("E29").Value = '  -3.06%  '
$ws.Range("E30").Value = '  +0.01%  '
$ws.Range("B31").Value = 'Cronos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.189'
$ws.Range("E31").Value = '  +2.59%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.135'
$ws.Range("E32").Value = '  +9.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '10.29'
$ws.Range("E33").Value = '  +11.03%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("E34").Value = '  +0.11%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '27.86'
$ws.Range("E35").Value = '  +2.14%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.150'
$ws.Range("E36").Value = '  -1.01%  '
$ws.Range("B37").Value = 'PolygonEcosystemToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.472'
$ws.Range("E37").Value = '  +5.37%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '7.17'
$ws.Range("E38").Value = '  -2.73%  '
$ws.Range("E39").Value = '  -0.80%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '24.80'
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '485.41'
$ws.Range("E41").Value = '  -5.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.65'
$ws.Range("E42").Value = '  +0.54%  '
$ws.Range("E43").Value = '  -3.01%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.769'
$ws.Range("E45").Value = '  -1.04%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.07'
$ws.Range("E46").Value = '  -6.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.94'
$ws.Range("E47").Value = '  -0.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '157.21'
$ws.Range("E48").Value = '  -2.51%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.847'
$ws.Range("E49").Value = '  +7.42%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.25'
$ws.Range("E51").Value = '  +4.13%  '
